# Weekly update: a new week's price record is inserted right after the
# existing row for 2022-01-04 (row 247), pushing every subsequent record
# down by one row. The new row is seeded as a copy of row 247 (so it
# inherits formatting/style) and then its date is corrected to the new
# week's date (2023-03-31, serial 45016).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(247).Copy()
$ws.Rows.Item(248).Insert()

$ws.Cells.Item(248, 4).Value = 45016
